# "Questions from James to consider"
# Insert a new row above the first timebox slot (old row 3) to hold an
# extra "Time Block" header/example entry, which pushes the existing
# schedule rows (old 3-20) down by one (new 4-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 3; Excel copies formatting from the
# row above (row 2) by default, which we'll override below for A3/B3.
$ws.Rows("3:3").Insert()

# A3: time value 6:40 AM, using the sheet's existing custom time format,
# but at the smaller (default) 12pt font rather than the big 18pt font
# used by the rest of column A.
$ws.Range("A3").NumberFormat = "[$-409]h:mm\ AM/PM;@"
$ws.Range("A3").Font.Size = 12
$ws.Range("A3").Value = 0.27777777777777779

# B3: new label "Time Block", left at the default (General) style / 12pt
# font instead of the big 18pt font used elsewhere in column B.
$ws.Range("B3").Font.Size = 12
$ws.Range("B3").Value = "Time Block"

# The trailing empty cell that rolled down into the new last row (C21)
# isn't present in the target sheet - drop it.
$ws.Range("C21").Clear()

# Move the active selection to B8.
[void]$ws.Range("B8").Select()
